$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns B, C, D, E, F, H, K, M across rows 2-25
# Using an array of row data to set values efficiently via Cells.Item

$ws.Cells.Item(2, 2).Value = 0.4365890793926042
$ws.Cells.Item(2, 3).Value = 0.02298481258995366
$ws.Cells.Item(2, 4).Value = 0.2179683731171878
$ws.Cells.Item(2, 5).Value = 0.1375878983431491
$ws.Cells.Item(2, 6).Value = 3.778795530854552
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 11).Value = 0.3818434170170804
$ws.Cells.Item(2, 13).Value = 0.2546306249913002

$ws.Cells.Item(3, 2).Value = 0.4149568748447621
$ws.Cells.Item(3, 3).Value = 0.02054100797720793
$ws.Cells.Item(3, 4).Value = 0.2057779472137469
$ws.Cells.Item(3, 5).Value = 0.1255670726398606
$ws.Cells.Item(3, 6).Value = 3.557895458523234
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 11).Value = 0.3602888624461258
$ws.Cells.Item(3, 13).Value = 0.2359885170281188

$ws.Cells.Item(4, 2).Value = 0.4021637306788364
$ws.Cells.Item(4, 3).Value = 0.01908389746481021
$ws.Cells.Item(4, 4).Value = 0.1982317192728829
$ws.Cells.Item(4, 5).Value = 0.1182492430580666
$ws.Cells.Item(4, 6).Value = 3.422424573983477
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 11).Value = 0.3475414940864994
$ws.Cells.Item(4, 13).Value = 0.2247481408713625

$ws.Cells.Item(5, 2).Value = 0.3970726832780826
$ws.Cells.Item(5, 3).Value = 0.0185008000308855
$ws.Cells.Item(5, 4).Value = 0.1951407180492737
$ws.Cells.Item(5, 5).Value = 0.115282654715628
$ws.Cells.Item(5, 6).Value = 3.367255022453151
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 11).Value = 0.3424684131767606
$ws.Cells.Item(5, 13).Value = 0.2202187676389755

$ws.Cells.Item(6, 2).Value = 0.3962346820408129
$ws.Cells.Item(6, 3).Value = 0.01840461606803956
$ws.Cells.Item(6, 4).Value = 0.1946264891468132
$ws.Cells.Item(6, 5).Value = 0.1147909799882783
$ws.Cells.Item(6, 6).Value = 3.358096186686737
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 11).Value = 0.3416333476346836
$ws.Cells.Item(6, 13).Value = 0.2194697414104283

$ws.Cells.Item(7, 2).Value = 0.4020945769551361
$ws.Cells.Item(7, 3).Value = 0.01907599063446952
$ws.Cells.Item(7, 4).Value = 0.1981900977192055
$ws.Cells.Item(7, 5).Value = 0.1182091723704488
$ws.Cells.Item(7, 6).Value = 3.421680399169134
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 11).Value = 0.3474725857254555
$ws.Cells.Item(7, 13).Value = 0.224686849777008

$ws.Cells.Item(8, 2).Value = 0.4290283276575906
$ws.Cells.Item(8, 3).Value = 0.02213304509355396
$ws.Cells.Item(8, 4).Value = 0.213777555029921
$ws.Cells.Item(8, 5).Value = 0.1334298213664411
$ws.Cells.Item(8, 6).Value = 3.702592361287145
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 11).Value = 0.374309707539453
$ws.Cells.Item(8, 13).Value = 0.2481597004513532

$ws.Cells.Item(9, 2).Value = 0.4857633054940038
$ws.Cells.Item(9, 3).Value = 0.02848278974884977
$ws.Cells.Item(9, 4).Value = 0.2438798400196021
$ws.Cells.Item(9, 5).Value = 0.1637951062884753
$ws.Cells.Item(9, 6).Value = 4.255001867998317
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 11).Value = 0.4308494698541949
$ws.Cells.Item(9, 13).Value = 0.2958547618691867

$ws.Cells.Item(10, 2).Value = 0.529892290527016
$ws.Cells.Item(10, 3).Value = 0.03337945442947898
$ws.Cells.Item(10, 4).Value = 0.2657439638990979
$ws.Cells.Item(10, 5).Value = 0.1864469207310293
$ws.Cells.Item(10, 6).Value = 4.662184838814369
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 11).Value = 0.4748443858160272
$ws.Cells.Item(10, 13).Value = 0.3319579783198847

$ws.Cells.Item(11, 2).Value = 0.550511754085079
$ws.Cells.Item(11, 3).Value = 0.03566069354941703
$ws.Cells.Item(11, 4).Value = 0.2756426529661553
$ws.Cells.Item(11, 5).Value = 0.1968321778963613
$ws.Cells.Item(11, 6).Value = 4.847798002161426
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 11).Value = 0.4954076795043534
$ws.Cells.Item(11, 13).Value = 0.348623269108387

$ws.Cells.Item(12, 2).Value = 0.5583991422751353
$ws.Cells.Item(12, 3).Value = 0.03653253243938082
$ws.Cells.Item(12, 4).Value = 0.279384733297519
$ws.Cells.Item(12, 5).Value = 0.2007768754031574
$ws.Cells.Item(12, 6).Value = 4.918146339120767
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 11).Value = 0.5032747244793541
$ws.Cells.Item(12, 13).Value = 0.3549695319699424

$ws.Cells.Item(13, 2).Value = 0.5566969134645774
$ws.Cells.Item(13, 3).Value = 0.03634440713214815
$ws.Cells.Item(13, 4).Value = 0.278579083577597
$ws.Cells.Item(13, 5).Value = 0.1999267728851279
$ws.Cells.Item(13, 6).Value = 4.902992779721728
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 11).Value = 0.5015768319548215
$ws.Cells.Item(13, 13).Value = 0.353601161263839

$ws.Cells.Item(14, 2).Value = 0.5511590609696952
$ws.Cells.Item(14, 3).Value = 0.03573225877359221
$ws.Cells.Item(14, 4).Value = 0.275950640438964
$ws.Cells.Item(14, 5).Value = 0.1971564671884636
$ws.Cells.Item(14, 6).Value = 4.853584354662019
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 11).Value = 0.4960532928093357
$ws.Cells.Item(14, 13).Value = 0.3491446647836014

$ws.Cells.Item(15, 2).Value = 0.5477773128519061
$ws.Cells.Item(15, 3).Value = 0.03535834759320267
$ws.Cells.Item(15, 4).Value = 0.2743398319633457
$ws.Cells.Item(15, 5).Value = 0.1954611531455583
$ws.Cells.Item(15, 6).Value = 4.823328332074823
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 11).Value = 0.4926804397611591
$ws.Cells.Item(15, 13).Value = 0.3464195738429936

$ws.Cells.Item(16, 2).Value = 0.5285558009487943
$ws.Cells.Item(16, 3).Value = 0.03323147307585828
$ws.Cells.Item(16, 4).Value = 0.2650961481469096
$ws.Cells.Item(16, 5).Value = 0.1857698805715913
$ws.Cells.Item(16, 6).Value = 4.650062769241231
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 11).Value = 0.4735116839949001
$ws.Cells.Item(16, 13).Value = 0.3308737879318286

$ws.Cells.Item(17, 2).Value = 0.5169042247267441
$ws.Cells.Item(17, 3).Value = 0.03194064625463966
$ws.Cells.Item(17, 4).Value = 0.2594136443578634
$ws.Cells.Item(17, 5).Value = 0.1798455863789741
$ws.Cells.Item(17, 6).Value = 4.543872093339985
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 11).Value = 0.4618938719457901
$ws.Cells.Item(17, 13).Value = 0.3213993759142681

$ws.Cells.Item(18, 2).Value = 0.5102537874047925
$ws.Cells.Item(18, 3).Value = 0.03120323769638844
$ws.Cells.Item(18, 4).Value = 0.2561407270969767
$ws.Cells.Item(18, 5).Value = 0.1764456730250075
$ws.Cells.Item(18, 6).Value = 4.482829796502614
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 11).Value = 0.4552632942940136
$ws.Cells.Item(18, 13).Value = 0.3159726564784719

$ws.Cells.Item(19, 2).Value = 0.5080108378278112
$ws.Cells.Item(19, 3).Value = 0.03095442180782015
$ws.Cells.Item(19, 4).Value = 0.255031787613035
$ws.Cells.Item(19, 5).Value = 0.1752958111805398
$ws.Cells.Item(19, 6).Value = 4.462167901098695
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 11).Value = 0.4530271372379389
$ws.Cells.Item(19, 13).Value = 0.31413914172321

$ws.Cells.Item(20, 2).Value = 0.5181392453341971
$ws.Cells.Item(20, 3).Value = 0.03207753340333852
$ws.Cells.Item(20, 4).Value = 0.2600190179680624
$ws.Cells.Item(20, 5).Value = 0.1804754498020884
$ws.Cells.Item(20, 6).Value = 4.55517253028242
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 11).Value = 0.4631252511343291
$ws.Cells.Item(20, 13).Value = 0.322405587529957

$ws.Cells.Item(21, 2).Value = 0.5527835047027168
$ws.Cells.Item(21, 3).Value = 0.03591184284947246
$ws.Cells.Item(21, 4).Value = 0.2767228460880062
$ws.Cells.Item(21, 5).Value = 0.1979698435059021
$ws.Cells.Item(21, 6).Value = 4.86809511623585
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 11).Value = 0.4976735053244568
$ws.Cells.Item(21, 13).Value = 0.3504526771804279

$ws.Cells.Item(22, 2).Value = 0.5758878521772317
$ws.Cells.Item(22, 3).Value = 0.03846444013440475
$ws.Cells.Item(22, 4).Value = 0.2876030434820791
$ws.Cells.Item(22, 5).Value = 0.2094737290249071
$ws.Cells.Item(22, 6).Value = 5.072964301113245
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 11).Value = 0.5207206222674472
$ws.Cells.Item(22, 13).Value = 0.3689902107869614

$ws.Cells.Item(23, 2).Value = 0.5635140291905429
$ws.Cells.Item(23, 3).Value = 0.03709771577386789
$ws.Cells.Item(23, 4).Value = 0.2817992779399958
$ws.Cells.Item(23, 5).Value = 0.2033273229359907
$ws.Cells.Item(23, 6).Value = 4.963587331918745
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 11).Value = 0.5083767616614523
$ws.Cells.Item(23, 13).Value = 0.3590771875058891

$ws.Cells.Item(24, 2).Value = 0.5175807424234335
$ws.Cells.Item(24, 3).Value = 0.03201563208870084
$ws.Cells.Item(24, 4).Value = 0.259745347190659
$ws.Cells.Item(24, 5).Value = 0.1801906697584883
$ws.Cells.Item(24, 6).Value = 4.550063577809425
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 11).Value = 0.4625683931307947
$ws.Cells.Item(24, 13).Value = 0.321950616249822

$ws.Cells.Item(25, 2).Value = 0.4699896866076756
$ws.Cells.Item(25, 3).Value = 0.02672546918509511
$ws.Cells.Item(25, 4).Value = 0.2357829018124562
$ws.Cells.Item(25, 5).Value = 0.1555223178531691
$ws.Cells.Item(25, 6).Value = 4.105355574994064
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 11).Value = 0.4151279556318457
$ws.Cells.Item(25, 13).Value = 0.2827689630747798
